# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45208 (2023-10-09) to serial date 45212 (2023-10-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 17; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45208) {
        $cell.Value = 45212
    }
}
